# Update automatico via Actualizar 04-15-2021 12-03-43
# Shifts the "Ultimo" (last-checked) timestamp column (D) forward for the
# three blocks of rows that share a common refresh timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-15: newest batch
$ws.Range("D2:D15").Value = 44301.50230610963

# Rows 16-29: previous batch, shifted down
$ws.Range("D16:D29").Value = 44301.48091100695

# Rows 30-43: oldest batch, shifted down
$ws.Range("D30:D43").Value = 44301.45951362269
